# Update LR-pair (Il15 -> Il2ra) NATMI edge-weight table with newly recomputed
# TPM-based ligand/receptor expression values.
#
# Column layout (row 1 header):
#   A Sending cluster          F Ligand detection rate     K Receptor-expressing cells   P Receptor derived specificity (total)
#   B Ligand symbol            G Ligand avg expr value     L Receptor detection rate     Q Edge avg expression weight
#   C Receptor symbol          H Ligand total expr value   M Receptor avg expr value     R Edge total expression weight
#   D Target cluster           I Ligand specificity (avg)  N Receptor total expr value   S Edge avg expression specificity
#   E Ligand-expressing cells  J Ligand specificity (tot)  O Receptor derived specificity(avg) T Edge total expression specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly recomputed per-cluster ligand (Il15) average expression values (TPM-based).
# The FAPs value happens to be numerically identical before/after the recompute.
$ligandAvg = @{
    "ECs"           = 6.654043666666666
    "FAPs"          = 1.987132666666667
    "MuSCs"         = 0.9593116666666667
    "Resolving-Mac" = 11.92023
}

# Newly recomputed per-cluster receptor (Il2ra) average expression values (TPM-based).
# The FAPs value happens to be numerically identical before/after the recompute.
$receptorAvg = @{
    "ECs"           = 0.2196916666666667
    "FAPs"          = 0.05768400000000001
    "MuSCs"         = 1.200477666666667
    "Resolving-Mac" = 0.534994
}

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# Ligand/receptor expressing-cell counts are unaffected by the TPM update, and are
# identical (3) for every cluster in this dataset, matching columns E and K*L.
$cellCount = 3

$ligandTotal    = @{}
$receptorTotal  = @{}
foreach ($c in $clusters) {
    $ligandTotal[$c]   = $ligandAvg[$c]   * $cellCount
    $receptorTotal[$c] = $receptorAvg[$c] * $cellCount
}

$sumLigandAvg   = 0
$sumReceptorAvg = 0
foreach ($c in $clusters) {
    $sumLigandAvg   += $ligandAvg[$c]
    $sumReceptorAvg += $receptorAvg[$c]
}

$ligandSpec   = @{}
$receptorSpec = @{}
foreach ($c in $clusters) {
    $ligandSpec[$c]   = $ligandAvg[$c]   / $sumLigandAvg
    $receptorSpec[$c] = $receptorAvg[$c] / $sumReceptorAvg
}

# The sheet lists, for every (Sending cluster, Target cluster) pair, a row, in a
# fixed block order (ECs, FAPs, MuSCs, Resolving-Mac as Sending cluster), each
# block containing the four Target clusters in the same order, starting at row 2.
#
# The raw average/total expression values (G,H for ligand; M,N for receptor) only
# actually change in the worksheet when the owning cluster's recomputed average
# differs from its old value -- true for every cluster except FAPs, whose TPM
# recompute happens to leave it numerically unchanged. The specificity columns
# (I,J,O,P) and all edge columns (Q,R,S,T) depend on sums/products across *all*
# clusters, so they change on every row regardless of whether that row's own
# cluster's raw value changed.
$row = 2
foreach ($send in $clusters) {
    foreach ($target in $clusters) {

        if ($send -ne "FAPs") {
            $ws.Cells.Item($row, 7).Value  = $ligandAvg[$send]     # G Ligand average expression value
            $ws.Cells.Item($row, 8).Value  = $ligandTotal[$send]   # H Ligand total expression value
        }
        $ws.Cells.Item($row, 9).Value  = $ligandSpec[$send]        # I Ligand derived specificity (avg)
        $ws.Cells.Item($row, 10).Value = $ligandSpec[$send]        # J Ligand derived specificity (total)

        if ($target -ne "FAPs") {
            $ws.Cells.Item($row, 13).Value = $receptorAvg[$target]   # M Receptor average expression value
            $ws.Cells.Item($row, 14).Value = $receptorTotal[$target] # N Receptor total expression value
        }
        $ws.Cells.Item($row, 15).Value = $receptorSpec[$target]      # O Receptor derived specificity (avg)
        $ws.Cells.Item($row, 16).Value = $receptorSpec[$target]      # P Receptor derived specificity (total)

        # Q/R are recomputed straight from the raw (non-normalized) averages/totals,
        # so when neither side's raw value actually moved (both Sending and Target
        # cluster are FAPs) the product is bit-for-bit identical and is left alone.
        if (-not ($send -eq "FAPs" -and $target -eq "FAPs")) {
            $ws.Cells.Item($row, 17).Value = $ligandAvg[$send]   * $receptorAvg[$target]      # Q Edge avg expression weight
            $ws.Cells.Item($row, 18).Value = $ligandTotal[$send] * $receptorTotal[$target]    # R Edge total expression weight
        }
        $ws.Cells.Item($row, 19).Value = $ligandSpec[$send]  * $receptorSpec[$target]     # S Edge avg expression specificity
        $ws.Cells.Item($row, 20).Value = $ligandSpec[$send]  * $receptorSpec[$target]     # T Edge total expression specificity

        $row++
    }
}
